$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'66.966.14"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.56%  "
$ws.Range("D3").Value = "'3.516.94"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.89%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'608.62"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.70%  "
$ws.Range("D6").Value = "'147.93"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.91%  "
$ws.Range("D7").Value = "'3.516.81"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.91%  "
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("D9").Value = "'0.477"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.68%  "
$ws.Range("E10").Value = "  -0.66%  "
$ws.Range("D11").Value = "'7.95"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.70%  "
$ws.Range("E12").Value = "  -1.63%  "
$ws.Range("D13").Value = "'0.0000218"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.58%  "
$ws.Range("D14").Value = "'4.111.63"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.07%  "
$ws.Range("D15").Value = "'31.95"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.21%  "
$ws.Range("D16").Value = "'3.515.38"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.58%  "
$ws.Range("D17").Value = "'67.024.46"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.33%  "
$ws.Range("E18").Value = "  -0.47%  "
$ws.Range("D19").Value = "'10.69"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +8.15%  "
$ws.Range("D20").Value = "'6.44"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.46%  "
$ws.Range("E21").Value = "  -0.24%  "
$ws.Range("D22").Value = "'437.85"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.60%  "
$ws.Range("E23").Value = "  -2.64%  "
$ws.Range("D24").Value = "'79.31"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.40%  "
$ws.Range("D25").Value = "'3.657.04"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.91%  "
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("E27").Value = "  -3.95%  "
$ws.Range("D28").Value = "'9.77"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.07%  "
$ws.Range("E29").Value = "  -4.34%  "
$ws.Range("E30").Value = "  +0.50%  "
$ws.Range("E31").Value = "  -3.54%  "
$ws.Range("E32").Value = "  -1.94%  "
$ws.Range("E33").Value = "  +0.06%  "
$ws.Range("D34").Value = "'25.45"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.37%  "
$ws.Range("D35").Value = "'3.508.55"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.80%  "
$ws.Range("E36").Value = "  -2.22%  "
$ws.Range("D37").Value = "'1.81"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.50%  "
$ws.Range("E38").Value = "  +0.53%  "
$ws.Range("E39").Value = "  +0.02%  "
$ws.Range("D40").Value = "'0.999"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.00%  "
$ws.Range("B41").Value = "Monero"
$ws.Range("C41").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D41").Value = "'173.21"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.33%  "
$ws.Range("B42").Value = "Hedera"
$ws.Range("C42").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D42").Value = "'0.0894"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.05%  "
$ws.Range("E43").Value = "  -0.03%  "
$ws.Range("E44").Value = "  -9.75%  "
$ws.Range("D45").Value = "'0.897"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.86%  "
$ws.Range("D46").Value = "'46.07"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.82%  "
$ws.Range("D47").Value = "'28.02"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.59%  "
$ws.Range("E48").Value = "  -0.21%  "
$ws.Range("E49").Value = "  -1.69%  "
$ws.Range("D50").Value = "'2.46"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.06%  "
$ws.Range("D51").Value = "'0.990"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.63%  "
